$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (shifts existing rows 14-24 down to 15-25)
$ws.Rows.Item(14).Insert()

# Populate the new row with the added entry (Liver / Liver Steatosis with Geographic Pattern)
$ws.Range("A14").Value2 = "Liver"
$ws.Range("B14").Value2 = "Liver Steatosis with Geographic Pattern "
$ws.Range("C14").Value2 = "Clip 1 B-mode"
$ws.Range("D14").Value2 = "https://youtu.be/m_H0po7LaIo"

# The row insert does not carry the hyperlink anchors down with it, so rebuild
# every hyperlink (including the new one) against their correct, shifted cells.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D3"), "https://youtu.be/zxTC0YBY2RY")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://youtu.be/K2Wbg7BgXy4")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://youtu.be/2kRZcpi70Aw")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://youtu.be/91M82AIMyu0")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/15o_Km86IzM")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/RhSUFLTmTl4")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://youtu.be/m_H0po7LaIo")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://youtu.be/DjI1kEnzfSQ")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/U3ydTsRwxok")
$ws.Hyperlinks.Add($ws.Range("D19"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D20"), "https://youtu.be/JvwODCASLYQ")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://youtu.be/pc-vbxSRTbs")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://youtu.be/Axbee4vjNtU")
$ws.Hyperlinks.Add($ws.Range("D23"), "https://youtu.be/qushjTAy6XQ")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/_FckFwJwynI")
$ws.Hyperlinks.Add($ws.Range("D25"), "https://youtu.be/z_oaRVxRz5s")

# Adding hyperlinks re-stamps the cell style; restore the plain hyperlink
# formatting used throughout column D (skipping the two "Coming soon" rows
# that never had a Variant/Detail or Link value).
$ws.Range("D3:D17").Style = "Collegamento ipertestuale"
$ws.Range("D19:D25").Style = "Collegamento ipertestuale"

# Refresh the worksheet's remembered sort range/condition to cover the new row,
# matching what Excel records after re-sorting the (already-ordered) table.
$sort = $ws.Sort
$sort.SetRange($ws.Range("A2:C23"))
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A23"))
$sort.Header = 0
$sort.Apply()

# Restore selection to the cell that was active after the edit
$ws.Range("D14").Select() | Out-Null
